$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 86:87. This leaves the current row 85 untouched in
# place, pushes the current row 86 (the "malla 70 unidades" entry) down to
# row 88, and opens up rows 86 and 87 as blank.
$ws.Range("86:87").Insert()

# The row that used to be at 85 (the "malla 50 unidades" entry) needs to be
# preserved at row 87 - copy it down before we overwrite row 85 below.
$ws.Rows.Item(85).Copy()
$ws.Rows.Item(87).PasteSpecial()

# Row 85: overwrite with the new "O'Higgins" / Primera record (2022-03-08).
$ws.Cells.Item(85, 1).Value = 11
$ws.Cells.Item(85, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(85, 3).Value = "Bíobío"
$ws.Cells.Item(85, 4).Value = 44628
$ws.Cells.Item(85, 5).Value = 8
$ws.Cells.Item(85, 6).Value = 100112024
$ws.Cells.Item(85, 7).Value = "Choclo"
$ws.Cells.Item(85, 8).Value = "Choclero"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 3000
$ws.Cells.Item(85, 11).Value = 200
$ws.Cells.Item(85, 12).Value = 200
$ws.Cells.Item(85, 13).Value = 200
$ws.Cells.Item(85, 14).Value = "$/unidad"
$ws.Cells.Item(85, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(85, 16).Value = 200
$ws.Cells.Item(85, 17).Value = 1
$ws.Cells.Item(85, 18).Value = "Hortaliza"

# Row 86: brand-new "O'Higgins" / Segunda record (2022-03-08).
$ws.Cells.Item(86, 1).Value = 11
$ws.Cells.Item(86, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(86, 3).Value = "Bíobío"
$ws.Cells.Item(86, 4).Value = 44628
$ws.Cells.Item(86, 5).Value = 8
$ws.Cells.Item(86, 6).Value = 100112024
$ws.Cells.Item(86, 7).Value = "Choclo"
$ws.Cells.Item(86, 8).Value = "Choclero"
$ws.Cells.Item(86, 9).Value = "Segunda"
$ws.Cells.Item(86, 10).Value = 3500
$ws.Cells.Item(86, 11).Value = 150
$ws.Cells.Item(86, 12).Value = 150
$ws.Cells.Item(86, 13).Value = 150
$ws.Cells.Item(86, 14).Value = "$/unidad"
$ws.Cells.Item(86, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(86, 16).Value = 150
$ws.Cells.Item(86, 17).Value = 1
$ws.Cells.Item(86, 18).Value = "Hortaliza"
